$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Day 2) tweaks ---
# Re-typed "Valid Anagram" now carries a trailing space (new shared string).
$ws.Range("C3").Value = "Valid Anagram "
# Status upgraded to a revision note.
$ws.Range("G3").Value = "S | R1: Q1"

# --- Row 5 (Day 4) tweak ---
$ws.Range("G5").Value = "S | R1: Q2"

# --- New Row 14 (Day 13) ---
$ws.Range("A14").Value = "Day 13"
# Copy the date-formatted style down from the row above, then overwrite the value.
$ws.Range("B13").Copy($ws.Range("B14"))
$ws.Range("B14").Value = 45815

$ws.Range("C14").Value = "Best Time to Buy and Sell Stock II.py"
$ws.Range("D14").Value = "Jump Game II.py"
$ws.Range("E14").Value = "Longest Common Prefix.py"
$ws.Range("F14").Value = "Greedy, Arrays, String"
$ws.Range("G14").Value = "S"

$h1 = $ws.Hyperlinks.Add($ws.Range("C14"), "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 13/Best Time to Buy and Sell Stock II.py", "", "Best Time to Buy and Sell Stock II.py")
$h1.TextToDisplay = "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 13/Best Time to Buy and Sell Stock II.py"
$ws.Range("C14").Value = "Best Time to Buy and Sell Stock II.py"

$h2 = $ws.Hyperlinks.Add($ws.Range("D14"), "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 13/Jump Game II.py", "", "Jump Game II.py")
$h2.TextToDisplay = "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 13/Jump Game II.py"
$ws.Range("D14").Value = "Jump Game II.py"

$h3 = $ws.Hyperlinks.Add($ws.Range("E14"), "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 13/Longest Common Prefix.py", "", "Longest Common Prefix.py")
$h3.TextToDisplay = "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/Day 13/Longest Common Prefix.py"
$ws.Range("E14").Value = "Longest Common Prefix.py"

# --- Selection moves on to the next empty row ---
$ws.Range("A15").Select()
